$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like plain numbers (e.g. "5.60"),
# which Excel would silently coerce to a numeric value (5.6) on assignment.
# Force text interpretation via NumberFormat "@" while writing, then restore
# the default "Normal" style so no stray style index is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.514.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.055.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.95%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.654"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +15.31%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.30%  "

$ws.Range("E10").Value = "  +2.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0771"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.38%  "

$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.919"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +26.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.356.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.042.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.387.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.58%  "

$ws.Range("E22").Value = "  +6.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.37%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.18%  "

$ws.Range("E26").Value = "  +4.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.118"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +23.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.35%  "

$ws.Range("E32").Value = "  +3.68%  "

$ws.Range("E33").Value = "  +9.49%  "

$ws.Range("E34").Value = "  +9.61%  "

$ws.Range("E35").Value = "  +6.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.10%  "

$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("E38").Value = "  +4.60%  "

$ws.Range("E39").Value = "  +17.06%  "

$ws.Range("E40").Value = "  +33.58%  "

$ws.Range("E41").Value = "  +18.47%  "

$ws.Range("E42").Value = "  +2.56%  "

$ws.Range("E43").Value = "  +4.53%  "

$ws.Range("E44").Value = "  +6.43%  "

$ws.Range("E45").Value = "  +4.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.39%  "

$ws.Range("E48").Value = "  +6.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.421.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.07%  "

$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
